$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 17:43"

# --- Update country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
# Columns: B=Casos totales C=Nuevos casos D=Casos activos E=Recuperados
#          F=Casos criticos G=Muertes hoy H=Muertes

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4652842
$ws.Range("C4").Value = 17857
$ws.Range("D4").Value = 2286492
$ws.Range("E4").Value = 2210634
$ws.Range("G4").Value = 431
$ws.Range("H4").Value = 155716

# Row 6 - India
$ws.Range("B6").Value = 1677853
$ws.Range("C6").Value = 38503
$ws.Range("D6").Value = 1081730
$ws.Range("E6").Value = 559938
$ws.Range("G6").Value = 399
$ws.Range("H6").Value = 36185

# Row 12 - España
$ws.Range("B12").Value = 335602
$ws.Range("C12").Value = 3092
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 28445

# Row 14 - Reino Unido
$ws.Range("B14").Value = 303181
$ws.Range("C14").Value = 880
$ws.Range("G14").Value = 120
$ws.Range("H14").Value = 46119

# Row 18 - Italia
$ws.Range("B18").Value = 247537
$ws.Range("C18").Value = 379
$ws.Range("E18").Value = 12600
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 35141

# Row 21 - Alemania
$ws.Range("B21").Value = 210112
$ws.Range("C21").Value = 459
$ws.Range("E21").Value = 8591

# Row 25 - Canada
$ws.Range("B25").Value = 115935
$ws.Range("C25").Value = 136
$ws.Range("D25").Value = 101030
$ws.Range("E25").Value = 5973
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 8932

# Rows 38/39 - Belgica & Republica Dominicana swap ranking order
# (Republica Dominicana overtakes Belgica with updated figures)
$ws.Range("A38").Value = "Republica Dominicana"
$ws.Range("B38").Value = 69649
$ws.Range("C38").Value = 1734
$ws.Range("D38").Value = 36470
$ws.Range("E38").Value = 32019
$ws.Range("G38").Value = 14
$ws.Range("H38").Value = 1160

$ws.Range("A39").Value = "Belgica"
$ws.Range("B39").Value = 68006
$ws.Range("C39").Value = 671
$ws.Range("D39").Value = 17513
$ws.Range("E39").Value = 40653
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 9840

# Row 45 - Singapur
$ws.Range("D45").Value = 46491
$ws.Range("E45").Value = 5687

# Row 85 - Senegal
$ws.Range("B85").Value = 10232
$ws.Range("C85").Value = 126
$ws.Range("D85").Value = 6776
$ws.Range("E85").Value = 3251
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 205

# Row 92 - Tayikistan
$ws.Range("B92").Value = 7409
$ws.Range("C92").Value = 43
$ws.Range("D92").Value = 6193
$ws.Range("E92").Value = 1156

# Row 141 - Jordania
$ws.Range("B141").Value = 1193
$ws.Range("C141").Value = 2
$ws.Range("D141").Value = 1084
$ws.Range("E141").Value = 98

# Row 159 - Reunion
$ws.Range("B159").Value = 660
$ws.Range("C159").Value = 3
$ws.Range("E159").Value = 64

# Row 177 - Islas Feroe
$ws.Range("D177").Value = 189
$ws.Range("E177").Value = 36
